$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 27225.25
$ws.Range("J9").Value = 2966.6667
$ws.Range("L9").Value = 2966.6667
$ws.Range("N9").Value = -3304.6667

$ws.Range("H12").Value = 393.6
$ws.Range("I12").Value = 398.44446
$ws.Range("K12").Value = 398.44446
$ws.Range("M12").Value = -228.44446

$ws.Range("H17").Value = 602463.75
$ws.Range("J17").Value = 602463.75
$ws.Range("L17").Value = 1807391.25
$ws.Range("N17").Value = -1807727.25

$ws.Range("H53").Value = 1136
$ws.Range("I53").Value = 155.16667
$ws.Range("J53").Value = 2116.8333
$ws.Range("K53").Value = 155.16667
$ws.Range("L53").Value = 2116.8333
$ws.Range("M53").Value = 481.83333
$ws.Range("N53").Value = -3390.8333

$ws.Range("H80").Value = 4772
$ws.Range("I80").Value = 4166.6665
$ws.Range("J80").Value = 4999
$ws.Range("K80").Value = 12499.9995
$ws.Range("L80").Value = 14997
$ws.Range("M80").Value = -11501.9995
$ws.Range("N80").Value = -16993

$ws.Range("H83").Value = 4772
$ws.Range("I83").Value = 4166.6665
$ws.Range("J83").Value = 4999
$ws.Range("K83").Value = 37499.9985
$ws.Range("L83").Value = 44991
$ws.Range("M83").Value = -32507.9985
$ws.Range("N83").Value = -54975

$ws.Range("H112").Value = 5749137
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 5749137
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 17247411
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -17249627

$ws.Range("H129").Value = 4398
$ws.Range("J129").Value = 5997
$ws.Range("L129").Value = 17991
$ws.Range("N129").Value = -27991

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1835.3
$ws.Range("I45").Value = 1589.25
$ws.Range("K45").Value = 1589.25
$ws.Range("M45").Value = -1212.25

$ws.Range("H61").Value = 11632516
$ws.Range("I61").Value = 12824826
$ws.Range("K61").Value = 12824826
$ws.Range("M61").Value = -12824614

$ws.Range("H102").Value = 3635.75
$ws.Range("I102").Value = 2300
$ws.Range("J102").Value = 4437.2
$ws.Range("K102").Value = 2300
$ws.Range("L102").Value = 4437.2
$ws.Range("M102").Value = -678
$ws.Range("N102").Value = -7681.2

$ws.Range("H132").Value = 27819474
$ws.Range("I132").Value = 2787.258
$ws.Range("J132").Value = 200282940
$ws.Range("K132").Value = 8361.773999999999
$ws.Range("L132").Value = 600848820
$ws.Range("M132").Value = -5831.773999999999
$ws.Range("N132").Value = -600853880

$ws.Range("H135").Value = 51483.855
$ws.Range("J135").Value = 55666.168
$ws.Range("L135").Value = 55666.168
$ws.Range("N135").Value = -65806.16800000001

$ws.Range("H136").Value = 11632516
$ws.Range("I136").Value = 12824826
$ws.Range("K136").Value = 38474478
$ws.Range("M136").Value = -38471928

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 949.5
$ws.Range("I22").Value = 499.5
$ws.Range("J22").Value = 1399.5
$ws.Range("K22").Value = 499.5
$ws.Range("L22").Value = 1399.5
$ws.Range("M22").Value = -326.5
$ws.Range("N22").Value = -1745.5

$ws.Range("H111").Value = 41444
$ws.Range("J111").Value = 41444
$ws.Range("L111").Value = 41444
$ws.Range("N111").Value = -49624

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2401.5293
$ws.Range("I58").Value = 1789.125
$ws.Range("K58").Value = 1789.125
$ws.Range("M58").Value = -1586.125

$ws.Range("H132").Value = 4035.7354
$ws.Range("I132").Value = 3487.1724
$ws.Range("K132").Value = 10461.5172
$ws.Range("M132").Value = -7931.5172

$ws.Range("H136").Value = 2401.5293
$ws.Range("I136").Value = 1789.125
$ws.Range("K136").Value = 5367.375
$ws.Range("M136").Value = -2817.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 47635292
$ws.Range("I4").Value = 63901468
$ws.Range("J4").Value = 27615384
$ws.Range("K4").Value = 191704404
$ws.Range("L4").Value = 82846152
$ws.Range("M4").Value = -191704292
$ws.Range("N4").Value = -82846376

$ws.Range("H22").Value = 127.5
$ws.Range("I22").Value = 124.333336
$ws.Range("J22").Value = 175
$ws.Range("K22").Value = 373.000008
$ws.Range("L22").Value = 525
$ws.Range("M22").Value = -204.000008
$ws.Range("N22").Value = -863

$ws.Range("H27").Value = 127.5
$ws.Range("I27").Value = 124.333336
$ws.Range("J27").Value = 175
$ws.Range("K27").Value = 373.000008
$ws.Range("L27").Value = 525
$ws.Range("M27").Value = -271.000008
$ws.Range("N27").Value = -729

$ws.Range("H107").Value = 1340
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 2150
$ws.Range("K107").Value = 2400
$ws.Range("L107").Value = 6450
$ws.Range("M107").Value = -480
$ws.Range("N107").Value = -10290

$ws.Range("H116").Value = 599.3333
$ws.Range("I116").Value = 399
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 1197
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 2245
$ws.Range("N116").Value = -9884

$ws.Range("H128").Value = 115925
$ws.Range("I128").Value = 115925
$ws.Range("K128").Value = 347775
$ws.Range("M128").Value = -342795

$ws.Range("H131").Value = 1376.25
$ws.Range("I131").Value = 728.3333
$ws.Range("J131").Value = 1765
$ws.Range("K131").Value = 2184.9999
$ws.Range("L131").Value = 5295
$ws.Range("M131").Value = 2855.0001
$ws.Range("N131").Value = -15375

$ws.Range("H133").Value = 4390.75
$ws.Range("I133").Value = 4390.75
$ws.Range("K133").Value = 13172.25
$ws.Range("M133").Value = -8112.25

$ws.Range("H134").Value = 10734.875
$ws.Range("J134").Value = 18127.25
$ws.Range("L134").Value = 54381.75
$ws.Range("N134").Value = -64521.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 30008
$ws.Range("I70").Value = 30008
$ws.Range("K70").Value = 30008
$ws.Range("M70").Value = -29738

$ws.Range("H73").Value = 30008
$ws.Range("I73").Value = 30008
$ws.Range("K73").Value = 30008
$ws.Range("M73").Value = -29072

$ws.Range("H112").Value = 55000
$ws.Range("I112").Value = 55000
$ws.Range("K112").Value = 55000
$ws.Range("M112").Value = -53892

$ws.Range("H122").Value = 1890.6923
$ws.Range("I122").Value = 1798.7778
$ws.Range("J122").Value = 2097.5
$ws.Range("K122").Value = 5396.3334
$ws.Range("L122").Value = 6292.5
$ws.Range("M122").Value = -2946.3334
$ws.Range("N122").Value = -11192.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4012.125
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4012.125
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 4012.125
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -4602.125

$ws.Range("H27").Value = 4012.125
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4012.125
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 4012.125
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -4226.125

$ws.Range("H46").Value = 1514.6154
$ws.Range("J46").Value = 3742.0908
$ws.Range("L46").Value = 3742.0908
$ws.Range("N46").Value = -4118.0908

$ws.Range("H122").Value = 3724.1304
$ws.Range("I122").Value = 2574.8333
$ws.Range("K122").Value = 7724.499899999999
$ws.Range("M122").Value = -5274.499899999999

$ws.Range("H132").Value = 1793.037
$ws.Range("I132").Value = 1776.48
$ws.Range("K132").Value = 5329.440000000001
$ws.Range("M132").Value = -2799.440000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5795.9165
$ws.Range("I62").Value = 3130.2
$ws.Range("J62").Value = 7700
$ws.Range("K62").Value = 3130.2
$ws.Range("L62").Value = 7700
$ws.Range("M62").Value = -2506.2
$ws.Range("N62").Value = -8948

$ws.Range("H65").Value = 5795.9165
$ws.Range("I65").Value = 3130.2
$ws.Range("J65").Value = 7700
$ws.Range("K65").Value = 15651
$ws.Range("L65").Value = 38500
$ws.Range("M65").Value = -12531
$ws.Range("N65").Value = -44740

$ws.Range("H122").Value = 2126.9
$ws.Range("I122").Value = 1914.6
$ws.Range("K122").Value = 5743.799999999999
$ws.Range("M122").Value = -3293.799999999999

$ws.Range("H126").Value = 6327.5386
$ws.Range("I126").Value = 6114.364
$ws.Range("K126").Value = 18343.092
$ws.Range("M126").Value = -15873.092
